$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 505.4
$ws.Cells.Item(8, 9).Value = 505.4
$ws.Cells.Item(8, 11).Value = 1516.2
$ws.Cells.Item(8, 13).Value = -1377.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 289.9375
$ws.Cells.Item(28, 9).Value = 257.42856
$ws.Cells.Item(28, 11).Value = 257.42856
$ws.Cells.Item(28, 13).Value = 227.57144

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3597
$ws.Cells.Item(76, 9).Value = 3246.25
$ws.Cells.Item(76, 11).Value = 3246.25
$ws.Cells.Item(76, 13).Value = -2931.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 3597
$ws.Cells.Item(79, 9).Value = 3246.25
$ws.Cells.Item(79, 11).Value = 3246.25
$ws.Cells.Item(79, 13).Value = -2154.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 66670920
$ws.Cells.Item(86, 9).Value = 100004200
$ws.Cells.Item(86, 10).Value = 33337644
$ws.Cells.Item(86, 11).Value = 100004200
$ws.Cells.Item(86, 12).Value = 33337644
$ws.Cells.Item(86, 13).Value = -100003077
$ws.Cells.Item(86, 14).Value = -33339890

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 66670920
$ws.Cells.Item(89, 9).Value = 100004200
$ws.Cells.Item(89, 10).Value = 33337644
$ws.Cells.Item(89, 11).Value = 500021000
$ws.Cells.Item(89, 12).Value = 166688220
$ws.Cells.Item(89, 13).Value = -500015384
$ws.Cells.Item(89, 14).Value = -166699452

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 4499.5
$ws.Cells.Item(98, 9).Value = 4499.5
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 4499.5
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).ClearContents()
$ws.Cells.Item(98, 14).Value = -3001.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2764.7273
$ws.Cells.Item(113, 10).Value = 2584.8
$ws.Cells.Item(113, 12).Value = 2584.8
$ws.Cells.Item(113, 14).Value = -9092.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 4499.5
$ws.Cells.Item(122, 9).Value = 4499.5
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 13498.5
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).Value = -11048.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2949554.8
$ws.Cells.Item(137, 9).Value = 5821.591
$ws.Cells.Item(137, 10).Value = 8346399
$ws.Cells.Item(137, 11).Value = 17464.773
$ws.Cells.Item(137, 12).Value = 25039197
$ws.Cells.Item(137, 13).Value = -14914.773
$ws.Cells.Item(137, 14).Value = -25044297

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 5122.891
$ws.Cells.Item(138, 9).Value = 6315.7036
$ws.Cells.Item(138, 10).Value = 3427.842
$ws.Cells.Item(138, 11).Value = 18947.1108
$ws.Cells.Item(138, 12).Value = 10283.526
$ws.Cells.Item(138, 13).Value = -13807.1108
$ws.Cells.Item(138, 14).Value = -20563.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1000
$ws.Cells.Item(2, 9).Value = 736.4545000000001
$ws.Cells.Item(2, 11).Value = 736.4545000000001
$ws.Cells.Item(2, 13).Value = -623.4545000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 56343.26
$ws.Cells.Item(45, 9).Value = 85020.664
$ws.Cells.Item(45, 10).Value = 7182
$ws.Cells.Item(45, 11).Value = 85020.664
$ws.Cells.Item(45, 12).Value = 7182
$ws.Cells.Item(45, 13).Value = -84643.664
$ws.Cells.Item(45, 14).Value = -7936

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 3468.3125
$ws.Cells.Item(110, 9).Value = 2044.2222
$ws.Cells.Item(110, 10).Value = 5299.2856
$ws.Cells.Item(110, 11).Value = 2044.2222
$ws.Cells.Item(110, 12).Value = 5299.2856
$ws.Cells.Item(110, 13).Value = 0.7778000000000702
$ws.Cells.Item(110, 14).Value = -9389.285599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1000
$ws.Cells.Item(116, 9).Value = 736.4545000000001
$ws.Cells.Item(116, 11).Value = 736.4545000000001
$ws.Cells.Item(116, 13).Value = 1557.5455

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 621.7143
$ws.Cells.Item(122, 9).Value = 379.6
$ws.Cells.Item(122, 10).Value = 1227
$ws.Cells.Item(122, 11).Value = 1138.8
$ws.Cells.Item(122, 12).Value = 3681
$ws.Cells.Item(122, 13).Value = 1311.2
$ws.Cells.Item(122, 14).Value = -8581

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2185.32
$ws.Cells.Item(132, 9).Value = 1486.25
$ws.Cells.Item(132, 11).Value = 4458.75
$ws.Cells.Item(132, 13).Value = -1928.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1000
$ws.Cells.Item(3, 9).Value = 736.4545000000001
$ws.Cells.Item(3, 11).Value = 736.4545000000001
$ws.Cells.Item(3, 13).Value = -622.4545000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 12809.723
$ws.Cells.Item(105, 9).Value = 10691.071
$ws.Cells.Item(105, 10).Value = 20225
$ws.Cells.Item(105, 11).Value = 10691.071
$ws.Cells.Item(105, 12).Value = 20225
$ws.Cells.Item(105, 13).Value = -8944.071
$ws.Cells.Item(105, 14).Value = -23719

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 17194.95
$ws.Cells.Item(107, 9).Value = 17836.79
$ws.Cells.Item(107, 10).Value = 5000
$ws.Cells.Item(107, 11).Value = 17836.79
$ws.Cells.Item(107, 12).Value = 5000
$ws.Cells.Item(107, 13).Value = -15916.79
$ws.Cells.Item(107, 14).Value = -8840

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 28126926
$ws.Cells.Item(134, 9).Value = 2014.1111
$ws.Cells.Item(134, 11).Value = 6042.3333
$ws.Cells.Item(134, 13).Value = -3507.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3220.4644
$ws.Cells.Item(31, 9).Value = 4015.6924
$ws.Cells.Item(31, 11).Value = 4015.6924
$ws.Cells.Item(31, 13).Value = -3720.6924

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3220.4644
$ws.Cells.Item(34, 9).Value = 4015.6924
$ws.Cells.Item(34, 11).Value = 4015.6924
$ws.Cells.Item(34, 13).Value = -3813.6924

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1821.8667
$ws.Cells.Item(107, 9).Value = 1818.091
$ws.Cells.Item(107, 11).Value = 1818.091
$ws.Cells.Item(107, 13).Value = 101.9090000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 308252.34
$ws.Cells.Item(4, 9).Value = 592.4706
$ws.Cells.Item(4, 11).Value = 1777.4118
$ws.Cells.Item(4, 13).Value = -1665.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 610.5714
$ws.Cells.Item(5, 9).Value = 610.5714
$ws.Cells.Item(5, 11).Value = 1831.7142
$ws.Cells.Item(5, 13).Value = -1719.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 4445946.5
$ws.Cells.Item(122, 9).Value = 5556108.5
$ws.Cells.Item(122, 10).Value = 5296.6665
$ws.Cells.Item(122, 11).Value = 50004976.5
$ws.Cells.Item(122, 12).Value = 47669.9985
$ws.Cells.Item(122, 13).Value = -50002526.5
$ws.Cells.Item(122, 14).Value = -52569.9985

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 610.5714
$ws.Cells.Item(135, 9).Value = 610.5714
$ws.Cells.Item(135, 11).Value = 5495.1426
$ws.Cells.Item(135, 13).Value = -2960.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 97166.5
$ws.Cells.Item(93, 10).Value = 95999.664
$ws.Cells.Item(93, 12).Value = 95999.664
$ws.Cells.Item(93, 14).Value = -99743.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 41667910
$ws.Cells.Item(102, 9).Value = 62501190
$ws.Cells.Item(102, 11).Value = 62501190
$ws.Cells.Item(102, 13).Value = -62499568

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 143949.28
$ws.Cells.Item(107, 10).Value = 1594.25
$ws.Cells.Item(107, 12).Value = 1594.25
$ws.Cells.Item(107, 14).Value = -5434.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 952420.75
$ws.Cells.Item(132, 9).Value = 2309.111
$ws.Cells.Item(132, 10).Value = 2852644
$ws.Cells.Item(132, 11).Value = 6927.333
$ws.Cells.Item(132, 12).Value = 8557932
$ws.Cells.Item(132, 13).Value = -4397.333
$ws.Cells.Item(132, 14).Value = -8562992

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2088.75
$ws.Cells.Item(40, 9).Value = 2167.5833
$ws.Cells.Item(40, 10).Value = 1852.25
$ws.Cells.Item(40, 11).Value = 2167.5833
$ws.Cells.Item(40, 12).Value = 1852.25
$ws.Cells.Item(40, 13).Value = -2031.5833
$ws.Cells.Item(40, 14).Value = -2124.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1156.6666
$ws.Cells.Item(93, 9).Value = 1195.2941
$ws.Cells.Item(93, 11).Value = 1195.2941
$ws.Cells.Item(93, 13).Value = 52.70589999999993

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3144
$ws.Cells.Item(122, 9).Value = 2859.8262
$ws.Cells.Item(122, 11).Value = 8579.4786
$ws.Cells.Item(122, 13).Value = -6129.4786

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3601.111
$ws.Cells.Item(132, 9).Value = 2936.6667
$ws.Cells.Item(132, 11).Value = 8810.000100000001
$ws.Cells.Item(132, 13).Value = -6280.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 500499.5
$ws.Cells.Item(81, 10).Value = 500499.5
$ws.Cells.Item(81, 12).Value = 1000999
$ws.Cells.Item(81, 14).Value = -1003121

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 500499.5
$ws.Cells.Item(84, 10).Value = 500499.5
$ws.Cells.Item(84, 12).Value = 5004995
$ws.Cells.Item(84, 14).Value = -5015603

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 3177317
$ws.Cells.Item(107, 9).Value = 2492.5
$ws.Cells.Item(107, 10).Value = 4084409.8
$ws.Cells.Item(107, 11).Value = 7477.5
$ws.Cells.Item(107, 12).Value = 12253229.4
$ws.Cells.Item(107, 13).Value = -5557.5
$ws.Cells.Item(107, 14).Value = -12257069.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2370.04
$ws.Cells.Item(132, 9).Value = 2046.2858
$ws.Cells.Item(132, 10).Value = 2782.0908
$ws.Cells.Item(132, 11).Value = 6138.857400000001
$ws.Cells.Item(132, 12).Value = 8346.2724
$ws.Cells.Item(132, 13).Value = -3608.857400000001
$ws.Cells.Item(132, 14).Value = -13406.2724
